$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column F ("Address"), shifting the existing District column to G
$ws.Columns("F:F").Insert()

$ws.Range("F2").Value = "Address"

$ws.Range("F4").Value = "Govt. High School Anandangar Hubli City"
$ws.Range("F5").Value = "G H S IngalagiKundagol"
$ws.Range("F6").Value = "G H S YarebudihalKundagol"
$ws.Range("F7").Value = "Govt. High School BammigattiKalaghatagi"
$ws.Range("F9").Value = "G H S HarlapuraKundagola"
$ws.Range("F11").Value = "Krishana High School HallikeriNavalgund"
$ws.Range("F12").Value = "G H S KoliwadHubi Rural"
$ws.Range("F13").Value = "S J R High School Noolvi Hubli"
$ws.Range("F14").Value = "H P S AladakattiKalaghatagi"
$ws.Range("F17").Value = "G H S YaliwalKungdol"
$ws.Range("F18").Value = "K L E Society H F Kattimani High SchoolHubli City"
$ws.Range("F19").Value = "Reshmi Almeelad Urdu High School Dayanand Colony Keshwapur Hubballi"
$ws.Range("F20").Value = "Adarsh Balika High School Ramnagar"
$ws.Range("F21").Value = "Smt Shantamma K Bhorashetti High School ShirakolNavalgund"
$ws.Range("F24").Value = "G H S PendargalliHubli City"
$ws.Range("F25").Value = "S G G High School GokulHubli"
$ws.Range("F26").Value = "G H S KusugalHubli"
$ws.Range("F27").Value = "AM Govt. Urdu High School Kundgol"
$ws.Range("F28").Value = "G H S Kalaghatagi"
$ws.Range("F29").Value = "Knupadanu High School Hubballi"
$ws.Range("F30").Value = "Govt. Urdu High School Kalaghatagi"
$ws.Range("F31").Value = "Govt. Adarsh Vidyalaya DastikoppaKalaghatagi"
$ws.Range("F32").Value = "Navalgund"
$ws.Range("F33").Value = "Model High School Navalagund"
$ws.Range("F34").Value = "G H S Navalur Dharwad City"
$ws.Range("F36").Value = "G H S Kusugal"
$ws.Range("F37").Value = "G H S HirenartiKundgol"
$ws.Range("F38").Value = "Govt. Adarsha Vidyalaya Dastikoppa Kalagahatagi"
$ws.Range("F39").Value = "Alm Govt. Urdu High School Sadashivnagar Old Hubli"
$ws.Range("F41").Value = "Dr. Ambedkar High School Katnur Hubli"
$ws.Range("F42").Value = "Shivappanna Jigalur High School Old Hubli"
$ws.Range("F44").Value = "G H S BagadageriKalaghatagi"
$ws.Range("F45").Value = "G H S HireharakuniKundgol"
$ws.Range("F46").Value = "G H S HallyalHubli Rural"
$ws.Range("F47").Value = "Lamington Girls High School Hubli City"
$ws.Range("F48").Value = "National Sports Resi High SchoolHubballi"
$ws.Range("F49").Value = "Govt. High School NavanagarHubli Rural"
$ws.Range("F50").Value = "Basel Mission Girls High School Karwar RoadHubli City"
$ws.Range("F51").Value = "Govt. Girls High School Kalghatagi"
$ws.Range("F52").Value = "Gurukul High School Dharwad City"
$ws.Range("F53").Value = "G H S SurashettikoppaKalaghatagi"
$ws.Range("F54").Value = "G H S ChalamattiKalaghatagi"
